# Applies the two "merge split runs back into one run" edits described
# by the diff:
#
#  1) Slide 1, subtitle shape, first paragraph:
#       "Выполнил" | ":" | " студент " | "3 " | "курса " | "Черных А. В."
#     becomes
#       "Выполнил" | ":" | " студент 3 курса Черных А. В."
#
#  2) Slide 8, body placeholder, "Среда разработки" paragraph:
#       "Среда разработки" | ": Visual Studio" | ", Rider, " | "GitHub"
#     becomes
#       "Среда разработки" | ": Visual Studio, Rider, GitHub"
#
# In both cases only the trailing runs are collapsed into a single run;
# the leading, untouched runs are left exactly as they were. Replacing
# the text of a precise Characters() sub-range (rather than the whole
# paragraph/shape) is what makes PowerPoint rewrite just that span as
# one run instead of doing a minimal prefix-preserving text patch.

$p = $ppt.ActivePresentation

# --- Edit 1: Slide 1 --------------------------------------------------
$slide1 = $p.Slides.Item(1)
$subtitle = $slide1.Shapes.Item(2)
$tr1 = $subtitle.TextFrame.TextRange
$para1 = $tr1.Paragraphs(1, 1)
# Characters(10, 29) spans " студент " + "3 " + "курса " + "Черных А. В."
# (runs 3-6 of paragraph 1), leaving "Выполнил" and ":" untouched.
$merge1 = $para1.Characters(10, 29)
$merge1.Text = " студент 3 курса Черных А. В."

# --- Edit 2: Slide 8 ----------------------------------------------------
$slide8 = $p.Slides.Item(8)
$body8 = $slide8.Shapes.Item(3)
$tr8 = $body8.TextFrame.TextRange
$para5 = $tr8.Paragraphs(5, 1)
# Characters(17, 30) spans ": Visual Studio" + ", Rider, " + "GitHub"
# (runs 2-4 of this paragraph), leaving "Среда разработки" untouched.
$merge2 = $para5.Characters(17, 30)
$merge2.Text = ": Visual Studio, Rider, GitHub"
